# Updated cryptos list on Sun Mar  3 04:12:24 UTC 2024 with GitHub Actions
#
# Columns: A=index(unchanged) B=Coin C=Link D=Price E=Volume(1h)
# D/E are stored as literal text in the workbook (e.g. "61.974.61",
# "  -0.52%  " with padding spaces). Some of the new Price strings
# (e.g. "130.40", "0.740", "43.88") are syntactically valid numbers, so a
# plain .Value assignment would make Excel coerce them into numeric cells
# and lose the trailing zero / exact text. Force each write through a
# "text" number format, then drop the format back to the default Normal
# style so no stray style survives the round trip.
function Set-TextCell {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple in-place Price / Volume(1h) refreshes -------------------------
Set-TextCell $ws "D2"  "62.043.11"
Set-TextCell $ws "E2"  "  -0.35%  "
Set-TextCell $ws "D3"  "3.435.48"
Set-TextCell $ws "E3"  "  -0.26%  "
Set-TextCell $ws "E4"  "  -0.15%  "
Set-TextCell $ws "D5"  "411.47"
Set-TextCell $ws "E5"  "  -0.47%  "
Set-TextCell $ws "D6"  "130.40"
Set-TextCell $ws "E6"  "  +0.31%  "
Set-TextCell $ws "D7"  "0.637"
Set-TextCell $ws "E7"  "  +1.26%  "
Set-TextCell $ws "E8"  "  +0.03%  "
Set-TextCell $ws "D9"  "0.740"
Set-TextCell $ws "E9"  "  -1.64%  "
Set-TextCell $ws "D10" "0.141"
Set-TextCell $ws "E10" "  +0.89%  "
Set-TextCell $ws "D11" "43.88"
Set-TextCell $ws "E11" "  +0.89%  "
Set-TextCell $ws "E12" "  +12.88%  "
Set-TextCell $ws "D13" "9.39"
Set-TextCell $ws "E13" "  +4.02%  "
Set-TextCell $ws "D14" "3.976.74"
Set-TextCell $ws "E14" "  -0.50%  "
Set-TextCell $ws "E15" "  +0.12%  "
Set-TextCell $ws "D16" "21.32"
Set-TextCell $ws "E16" "  +3.06%  "
Set-TextCell $ws "D17" "3.430.19"
Set-TextCell $ws "E17" "  -0.69%  "
Set-TextCell $ws "D18" "12.52"
Set-TextCell $ws "E18" "  -5.29%  "
Set-TextCell $ws "E19" "  +1.84%  "
Set-TextCell $ws "D20" "61.988.73"
Set-TextCell $ws "E20" "  -0.35%  "
Set-TextCell $ws "D21" "508.61"
Set-TextCell $ws "E21" "  +27.54%  "
Set-TextCell $ws "D22" "93.19"
Set-TextCell $ws "E22" "  +3.08%  "
Set-TextCell $ws "D23" "3.32"
Set-TextCell $ws "E23" "  +3.46%  "
Set-TextCell $ws "D24" "13.53"
Set-TextCell $ws "E24" "  +0.89%  "
Set-TextCell $ws "D25" "3.36"
Set-TextCell $ws "E25" "  +1.90%  "
Set-TextCell $ws "D26" "35.08"
Set-TextCell $ws "E26" "  +2.76%  "
Set-TextCell $ws "D27" "9.27"
Set-TextCell $ws "E27" "  +6.39%  "
Set-TextCell $ws "E28" "  -0.04%  "
Set-TextCell $ws "D29" "12.19"
Set-TextCell $ws "E29" "  +2.13%  "
Set-TextCell $ws "E30" "  -1.81%  "
Set-TextCell $ws "E31" "  -2.02%  "
Set-TextCell $ws "E32" "  -2.53%  "
Set-TextCell $ws "D33" "42.08"
Set-TextCell $ws "E33" "  -4.69%  "

# --- Row 34 / 35: Dai and OKB swap ranking order, with refreshed values ---
Set-TextCell $ws "B34" "OKB"
Set-TextCell $ws "C34" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D34" "59.71"
Set-TextCell $ws "E34" "  +13.70%  "

Set-TextCell $ws "B35" "Dai"
Set-TextCell $ws "C35" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D35" "1.00"
Set-TextCell $ws "E35" "  +0.03%  "

Set-TextCell $ws "D36" "0.0502"
Set-TextCell $ws "E36" "  +0.48%  "

# --- Rows 37 / 38 / 39: FirstDigitalUSD / Stellar / LidoDAOToken rotate ---
Set-TextCell $ws "B37" "LidoDAOToken"
Set-TextCell $ws "C37" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D37" "3.49"
Set-TextCell $ws "E37" "  +2.62%  "

Set-TextCell $ws "B38" "FirstDigitalUSD"
Set-TextCell $ws "C38" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws "D38" "0.998"
Set-TextCell $ws "E38" "  -0.04%  "

Set-TextCell $ws "B39" "Stellar"
Set-TextCell $ws "C39" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D39" "0.139"
Set-TextCell $ws "E39" "  +5.78%  "

Set-TextCell $ws "D40" "2.75"
Set-TextCell $ws "E40" "  +17.34%  "
Set-TextCell $ws "D41" "148.41"
Set-TextCell $ws "E41" "  +5.50%  "

# --- Rows 42 / 43: Stacks and ARBITRUM swap ranking order -----------------
Set-TextCell $ws "B42" "ARBITRUM"
Set-TextCell $ws "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D42" "2.14"
Set-TextCell $ws "E42" "  +6.96%  "

Set-TextCell $ws "B43" "Stacks"
Set-TextCell $ws "C43" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws "D43" "2.95"
Set-TextCell $ws "E43" "  +1.58%  "

Set-TextCell $ws "D44" "0.320"
Set-TextCell $ws "E44" "  +1.27%  "
Set-TextCell $ws "D45" "4.35"
Set-TextCell $ws "E45" "  +6.73%  "
Set-TextCell $ws "E46" "  +22.69%  "
Set-TextCell $ws "D47" "16.69"
Set-TextCell $ws "E47" "  -1.27%  "
Set-TextCell $ws "D48" "121.41"
Set-TextCell $ws "E48" "  +22.56%  "
Set-TextCell $ws "D49" "23.06"
Set-TextCell $ws "E49" "  +0.66%  "
Set-TextCell $ws "D50" "0.147"
Set-TextCell $ws "E50" "  +19.55%  "
Set-TextCell $ws "D51" "2.146.51"
Set-TextCell $ws "E51" "  +0.77%  "
